$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths (A:D = 30 "characters") ---
# The engine stores ColumnWidth with a constant +5/6 offset versus the
# resulting OOXML <col width="..."> value, so compensate for it here in
# order to land on an exact width of 30 in the saved file.
$colOffset = 5 / 6
$ws.Columns.Item(1).ColumnWidth = 30 - $colOffset
$ws.Columns.Item(2).ColumnWidth = 30 - $colOffset
$ws.Columns.Item(3).ColumnWidth = 30 - $colOffset
$ws.Columns.Item(4).ColumnWidth = 30 - $colOffset

# --- Header row ---
$ws.Range("A1").Value = "nome"
$ws.Range("B1").Value = "peso"
$ws.Range("C1").Value = "zone"
$ws.Range("D1").Value = "partenza"

# --- Data rows ---
$ws.Range("A2").Value = "40032 (interno)"
$ws.Range("B2").Value = "CAMPO VUOTO"
$ws.Range("C2").Value = "CAMPO VUOTO"
$ws.Range("D2").Value = 0

$ws.Range("A3").Value = "40032 (interno)"
$ws.Range("B3").Value = "CAMPO VUOTO"
$ws.Range("C3").Value = "CAMPO VUOTO"
$ws.Range("D3").Value = 0

# --- Yellow fill on B2:C3 ("CAMPO VUOTO" cells) ---
# Build the fill on a scratch cell far away first (so both the foreground
# AND background color land on the finished style in one go), then copy
# that formatting onto the target cells and wipe the scratch cell again.
$yellow = 65535   # RGB(255,255,0) packed as BGR COM color value
$helper = $ws.Range("Z100")
$helper.Interior.Color = $yellow
$helper.Interior.PatternColor = $yellow

$helper.Copy()
$ws.Range("B2:C3").PasteSpecial(-4122)
$helper.Clear()

$excel.CutCopyMode = $false
